# Update the EffectiveDate column (F) for both data rows from 07302023 to 08302023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "08302023"
$ws.Range("F3").Value = "08302023"

# Update the selected cell in the sheet view to E7
$ws.Range("E7").Select()
